$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.851.29"
$ws.Range("E2").Value = "  +0.08%  "

$ws.Range("D3").Value = "2.538.10"
$ws.Range("E3").Value = "  -0.10%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.76%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.64"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.18%  "

$ws.Range("E7").Value = "  +0.60%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.546"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.30%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.94"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.03%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0824"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.38%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.77"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.93%  "

$ws.Range("E13").Value = "  -1.41%  "

$ws.Range("D14").Value = "2.929.78"
$ws.Range("E14").Value = "  +0.14%  "

$ws.Range("D15").Value = "2.552.21"
$ws.Range("E15").Value = "  -1.47%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.17"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.99%  "

$ws.Range("E17").Value = "  -0.26%  "

$ws.Range("D18").Value = "42.872.51"
$ws.Range("E18").Value = "  +0.22%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.05%  "

$ws.Range("D20").Value = "0.0₃0990"
$ws.Range("E20").Value = "  +1.02%  "

$ws.Range("E21").Value = "  +0.50%  "

$ws.Range("E22").Value = "  -0.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "254.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.32%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.97%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.89%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "27.78"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.43%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.43%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.32"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +9.82%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.14%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.91"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.21"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.49%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "156.97"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.79%  "

$ws.Range("E33").Value = "  +0.08%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.21"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +10.08%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0800"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.82%  "

$ws.Range("E36").Value = "  -2.16%  "

$ws.Range("E37").Value = "  -4.71%  "

$ws.Range("E38").Value = "  +1.26%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "24.63"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.29%  "

$ws.Range("E40").Value = "  +0.68%  "

$ws.Range("E41").Value = "  +9.98%  "

$ws.Range("E42").Value = "  +1.26%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.90"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.83%  "

$ws.Range("D44").Value = "2.083.23"
$ws.Range("E44").Value = "  -0.16%  "

$ws.Range("E45").Value = "  -1.53%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "86.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.75%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.99"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.20%  "

$ws.Range("D49").Value = "2.786.46"
$ws.Range("E49").Value = "  +0.07%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.64%  "

$ws.Range("E51").Value = "  +1.89%  "
